# Updated cryptos list on Fri Nov 10 06:00:08 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Simple per-row Price(D) / Volume(E) updates ---
# Each entry: row, newD (or $null to leave D unchanged), newE
$rows = @(
    @(2,  "36.736.86",   "  +0.30%  "),
    @(3,  "2.124.59",    "  +10.67%  "),
    @(4,  $null,         "  -0.01%  "),
    @(5,  "256.06",      "  +2.81%  "),
    @(6,  "0.669",       "  -4.38%  "),
    @(7,  $null,         "  +0.04%  "),
    @(8,  "47.23",       "  +6.10%  "),
    @(9,  "59.79",       "  +1.53%  "),
    @(10, "0.375",       "  +2.27%  "),
    @(11, $null,         "  -2.36%  "),
    @(12, $null,         "  +0.45%  "),
    @(13, "2.436.18",    "  +10.76%  "),
    @(14, "14.38",       "  -1.37%  "),
    @(15, $null,         "  +5.18%  "),
    @(16, "2.124.10",    "  +10.66%  "),
    @(17, "5.15",        "  +0.47%  "),
    @(18, "36.710.05",   "  +0.06%  "),
    @(19, "73.87",       "  -0.41%  "),
    @(20, "0.0₃0840",    "  -2.39%  "),
    @(21, $null,         "  +0.39%  "),
    @(22, "242.11",      "  -3.83%  "),
    @(23, "5.22",        "  +0.30%  "),
    @(24, $null,         "  +0.05%  "),
    @(25, $null,         "  -7.14%  "),
    @(28, "9.27",        "  +5.10%  "),
    @(29, $null,         "  -7.35%  "),
    @(30, "29.66",       "  +65.35%  "),
    @(31, $null,         "  -4.22%  "),
    @(32, $null,         "  -0.52%  "),
    @(33, "0.0963",      "  +13.53%  "),
    @(34, "0.0602",      "  -3.39%  "),
    @(35, $null,         "  +18.34%  "),
    @(36, "0.966",       "  +10.59%  "),
    @(37, $null,         "  -4.83%  "),
    @(38, $null,         "  +0.00%  "),
    @(39, $null,         "  -4.40%  "),
    @(40, $null,         "  -9.61%  "),
    @(41, $null,         "  +8.24%  "),
    @(42, "0.0227",      "  -0.73%  "),
    @(43, "99.49",       "  -6.26%  "),
    @(44, $null,         "  +10.03%  "),
    @(45, "16.28",       "  -5.56%  "),
    @(46, "1.361.07",    "  +1.73%  "),
    @(47, $null,         "  +12.61%  "),
    @(48, "0.0844",      "  +3.49%  "),
    @(49, "2.322.71",    "  +10.67%  ")
)

foreach ($r in $rows) {
    $rowNum = $r[0]
    $newD = $r[1]
    $newE = $r[2]
    if ($null -ne $newD) {
        $cellD = $ws.Cells.Item($rowNum, 4)
        $cellD.NumberFormat = "@"
        $cellD.Value = $newD
    }
    $ws.Cells.Item($rowNum, 5).Value = $newE
}

# --- Row 26/27 swap: Monero <-> EthereumClassic ---
$ws.Cells.Item(26, 2).Value = "EthereumClassic"
$ws.Cells.Item(26, 3).Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$cellD26 = $ws.Cells.Item(26, 4)
$cellD26.NumberFormat = "@"
$cellD26.Value = "21.96"
$ws.Cells.Item(26, 5).Value = "  +16.95%  "

$ws.Cells.Item(27, 2).Value = "Monero"
$ws.Cells.Item(27, 3).Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$cellD27 = $ws.Cells.Item(27, 4)
$cellD27.NumberFormat = "@"
$cellD27.Value = "172.22"
$ws.Cells.Item(27, 5).Value = "  +2.56%  "

# --- Row 50/51 swap: MXToken <-> RenderToken ---
$ws.Cells.Item(50, 2).Value = "RenderToken"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$cellD50 = $ws.Cells.Item(50, 4)
$cellD50.NumberFormat = "@"
$cellD50.Value = "2.30"
$ws.Cells.Item(50, 5).Value = "  -3.10%  "

$ws.Cells.Item(51, 2).Value = "MXToken"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$cellD51 = $ws.Cells.Item(51, 4)
$cellD51.NumberFormat = "@"
$cellD51.Value = "2.84"
$ws.Cells.Item(51, 5).Value = "  +1.86%  "
